# Generate Report for Handoff
# Updates the status of the bf6e3734-... file (row 3 on each sheet) to
# "Ready for handoff" and records the new handoff datetimes for the
# zh-cn and de-de target languages.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status columns (zh-cn = B, de-de = C) for row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: Status (B3) and Latest Handoff Datetime (D3) for row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-10 03:23:25"

# --- de-de sheet: Status (B3) and Latest Handoff Datetime (D3) for row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-10 03:23:28"
